# Auto-generated corrections to H:N profit-calc columns across all 8 sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1166002.6
$ws.Range("J17").Value = 1196675.2
$ws.Range("L17").Value = 3590025.6
$ws.Range("N17").Value = -3590361.6
$ws.Range("H33").Value = 750.2941
$ws.Range("I33").Value = 236.35715
$ws.Range("K33").Value = 236.35715
$ws.Range("M33").Value = -7.35714999999999
$ws.Range("H70").Value = 12505045
$ws.Range("I70").Value = 28576386
$ws.Range("K70").Value = 85729158
$ws.Range("M70").Value = -85728888
$ws.Range("H73").Value = 12505045
$ws.Range("I73").Value = 28576386
$ws.Range("K73").Value = 85729158
$ws.Range("M73").Value = -85728222
$ws.Range("H107").Value = 5431.7188
$ws.Range("I107").Value = 6224.0527
$ws.Range("J107").Value = 4273.6924
$ws.Range("K107").Value = 6224.0527
$ws.Range("L107").Value = 4273.6924
$ws.Range("M107").Value = -4304.0527
$ws.Range("N107").Value = -8113.6924
$ws.Range("H116").Value = 622860
$ws.Range("I116").Value = 1393470.6
$ws.Range("J116").Value = 6371.5
$ws.Range("K116").Value = 1393470.6
$ws.Range("L116").Value = 6371.5
$ws.Range("M116").Value = -1390028.6
$ws.Range("N116").Value = -13255.5
$ws.Range("H137").Value = 355115.47
$ws.Range("I137").Value = 526589.4
$ws.Range("K137").Value = 1579768.2
$ws.Range("M137").Value = -1577218.2
$ws.Range("H138").Value = 5033.75
$ws.Range("I138").Value = 527.82355
$ws.Range("J138").Value = 6663.553
$ws.Range("K138").Value = 1583.47065
$ws.Range("L138").Value = 19990.659
$ws.Range("M138").Value = 3556.52935
$ws.Range("N138").Value = -30270.659
$ws.Range("H141").Value = 5450.609
$ws.Range("J141").Value = 5742.222
$ws.Range("L141").Value = 17226.666
$ws.Range("N141").Value = -27586.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 520.44446
$ws.Range("I5").Value = 280.66666
$ws.Range("K5").Value = 280.66666
$ws.Range("M5").Value = -168.66666
$ws.Range("H32").Value = 3103.75
$ws.Range("I32").Value = 3080.4412
$ws.Range("K32").Value = 3080.4412
$ws.Range("M32").Value = -2793.4412
$ws.Range("H41").Value = 9842.5
$ws.Range("I41").Value = 7811
$ws.Range("K41").Value = 7811
$ws.Range("M41").Value = -7397
$ws.Range("H63").Value = 1956.2858
$ws.Range("I63").Value = 1938.8
$ws.Range("K63").Value = 1938.8
$ws.Range("M63").Value = -1252.8
$ws.Range("H66").Value = 1956.2858
$ws.Range("I66").Value = 1938.8
$ws.Range("K66").Value = 9694
$ws.Range("M66").Value = -6262
$ws.Range("H97").Value = 19482.37
$ws.Range("I97").Value = 9239.857
$ws.Range("J97").Value = 48161.4
$ws.Range("K97").Value = 9239.857
$ws.Range("L97").Value = 48161.4
$ws.Range("M97").Value = -8743.857
$ws.Range("N97").Value = -49153.4
$ws.Range("H102").Value = 3331.4375
$ws.Range("I102").Value = 3238.5
$ws.Range("K102").Value = 3238.5
$ws.Range("M102").Value = -1616.5
$ws.Range("H122").Value = 3148266
$ws.Range("I122").Value = 4484.5
$ws.Range("J122").Value = 7339974.5
$ws.Range("K122").Value = 13453.5
$ws.Range("L122").Value = 22019923.5
$ws.Range("M122").Value = -11003.5
$ws.Range("N122").Value = -22024823.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 520.44446
$ws.Range("I4").Value = 280.66666
$ws.Range("K4").Value = 280.66666
$ws.Range("M4").Value = -165.66666
$ws.Range("H94").Value = 3057.3928
$ws.Range("I94").Value = 3046.611
$ws.Range("J94").Value = 3076.8
$ws.Range("K94").Value = 3046.611
$ws.Range("L94").Value = 3076.8
$ws.Range("M94").Value = -2595.611
$ws.Range("N94").Value = -3978.8
$ws.Range("H99").Value = 52933.2
$ws.Range("I99").Value = 64916.5
$ws.Range("K99").Value = 64916.5
$ws.Range("M99").Value = -63418.5
$ws.Range("H134").Value = 2087.139
$ws.Range("I134").Value = 1892.2646
$ws.Range("K134").Value = 5676.793799999999
$ws.Range("M134").Value = -3141.793799999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3729.7727
$ws.Range("I31").Value = 2709.1177
$ws.Range("K31").Value = 2709.1177
$ws.Range("M31").Value = -2414.1177
$ws.Range("H34").Value = 3729.7727
$ws.Range("I34").Value = 2709.1177
$ws.Range("K34").Value = 2709.1177
$ws.Range("M34").Value = -2507.1177
$ws.Range("H58").Value = 2620.5
$ws.Range("I58").Value = 1801.0625
$ws.Range("J58").Value = 3348.889
$ws.Range("K58").Value = 1801.0625
$ws.Range("L58").Value = 3348.889
$ws.Range("M58").Value = -1598.0625
$ws.Range("N58").Value = -3754.889
$ws.Range("H62").Value = 94142.71000000001
$ws.Range("I62").Value = 4000
$ws.Range("K62").Value = 4000
$ws.Range("M62").Value = -3376
$ws.Range("H65").Value = 94142.71000000001
$ws.Range("I65").Value = 4000
$ws.Range("K65").Value = 20000
$ws.Range("M65").Value = -16880
$ws.Range("H134").Value = 2089338.9
$ws.Range("I134").Value = 2610342.5
$ws.Range("K134").Value = 7831027.5
$ws.Range("M134").Value = -7828492.5
$ws.Range("H136").Value = 2620.5
$ws.Range("I136").Value = 1801.0625
$ws.Range("J136").Value = 3348.889
$ws.Range("K136").Value = 5403.1875
$ws.Range("L136").Value = 10046.667
$ws.Range("M136").Value = -2853.1875
$ws.Range("N136").Value = -15146.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2506
$ws.Range("I136").Value = 2506
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7518
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = -2418

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 13579.889
$ws.Range("I122").Value = 12952.308
$ws.Range("K122").Value = 38856.924
$ws.Range("M122").Value = -36406.924
$ws.Range("H123").Value = 17368.37
$ws.Range("J123").Value = 17368.37
$ws.Range("L123").Value = 17368.37
$ws.Range("N123").Value = -22268.37
$ws.Range("H141").Value = 100429
$ws.Range("J141").Value = 100429
$ws.Range("L141").Value = 100429
$ws.Range("N141").Value = -110789

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 42804
$ws.Range("I38").Value = 14030
$ws.Range("K38").Value = 14030
$ws.Range("M38").Value = -13620

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 33752.25
$ws.Range("I29").Value = 40336.332
$ws.Range("J29").Value = 14000
$ws.Range("K29").Value = 40336.332
$ws.Range("L29").Value = 14000
$ws.Range("M29").Value = -40046.332
$ws.Range("N29").Value = -14580
$ws.Range("H122").Value = 3860.6206
$ws.Range("I122").Value = 2725.5
$ws.Range("K122").Value = 8176.5
$ws.Range("M122").Value = -5726.5
$ws.Range("H141").Value = 70715
$ws.Range("J141").Value = 70715
$ws.Range("L141").Value = 70715
$ws.Range("N141").Value = -81075

# Cells removed by this update (no longer applicable)
$wb.Worksheets.Item("CUL").Range("M136").ClearContents()

Write-Host "Applied Siren_Profits corrections"
